$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 9 ("9. Feature to download/upload list of words") row 10:
#   Status: In Progress -> Upload feature remained ; Assignee stays Arthur
$ws.Range("C10").Value = "Upload feature remained"
$ws.Range("D10").Value = "Arthur"

# Task 1 ("1. Reduce count of Gomel-sat news to 45 (3 pages)") row 2:
#   Status: Open -> In progress ; Assignee: (blank) -> Arthur
$ws.Range("C2").Value = "In progress"
$ws.Range("D2").Value = "Arthur"

# Widen the Status column (C) to fit the longer status text
$ws.Columns.Item(3).ColumnWidth = 35.8

# Move the active selection to C6
$ws.Range("C6").Select() | Out-Null
